$wb = $excel.ActiveWorkbook

# --- Rename sheets (task-order ids refreshed) ---
$wb.Worksheets.Item(1).Name = "GNG_TO-1651168698087312"
$wb.Worksheets.Item(2).Name = "NB_TO-1651168701107487"
$wb.Worksheets.Item(3).Name = "RS_TO-16511687011094782"
$wb.Worksheets.Item(4).Name = "TOL_TO-16511687011714807"
$wb.Worksheets.Item(5).Name = "vSAT_TO-1651168701249513"

# --- Sheet 1: GNG_TO ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("B2").Value = "go_stims-16511686980563128.csv"
$ws1.Range("B3").Value = "GNG_stims-1651168698071312.csv"
$ws1.Range("B4").Value = "go_stims-16511686980723164.csv"
$ws1.Range("B5").Value = "GNG_stims-16511686980863106.csv"

# --- Sheet 2: NB_TO ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("B2").Value = "TB-16511687010074766.csv"
$ws2.Range("B3").Value = "ZB-match_0-16511686986763206.csv"
$ws2.Range("B4").Value = "OB-16511687001125364.csv"
$ws2.Range("B5").Value = "TB-165116870030151.csv"
$ws2.Range("B6").Value = "ZB-match_2-16511686988243175.csv"
$ws2.Range("B7").Value = "ZB-match_6-16511686982903147.csv"
$ws2.Range("B8").Value = "OB-1651168699644314.csv"
$ws2.Range("B9").Value = "OB-1651168700138498.csv"
$ws2.Range("B10").Value = "TB-16511687010904787.csv"

# --- Sheet 3: RS_TO (name change only, no data changes) ---

# --- Sheet 4: TOL_TO ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("B2").Value = "MM_stims-16511687011234758.csv"
$ws4.Range("B3").Value = "ZM_stims-1651168701111477.csv"
$ws4.Range("B4").Value = "MM_stims-1651168701154512.csv"
$ws4.Range("B5").Value = "ZM_stims-16511687011234758.csv"
$ws4.Range("B6").Value = "MM_stims-16511687011704917.csv"
$ws4.Range("B7").Value = "ZM_stims-16511687011554773.csv"

# --- Sheet 5: vSAT_TO ---
$ws5 = $wb.Worksheets.Item(5)
$ws5.Range("B2").Value = "vSAT_stims-16511687012334793.csv"
$ws5.Range("B3").Value = "vSAT_stims-1651168701217489.csv"
$ws5.Range("B4").Value = "SAT_stims-16511687011774805.csv"
$ws5.Range("B5").Value = "SAT_stims-16511687012034788.csv"
